$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Truncate the tail of the existing run so it ends with "...which would "
#    (dropping "subsequently lower the annual income by a large margin. ")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "hire more people to do it which would subsequently lower the annual income by a large margin. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "hire more people to do it which would ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Append the whole new block of text right after "which would " (still
#    plain, no superscripts yet -- those get applied in step 3).
# ---------------------------------------------------------------------------
$newBlock = "mean a huge increase in salary and wages expense, increase salary and wages payables and subsequently reduce fiscal net income by a large margin and no one wants that. I mean, sure that" + [char]8217 + "ll be great for  our tax return mostly because we won" + [char]8217 + "t have to pay as much but It wouldn" + [char]8217 + "t bode well for future potential investors when they read our financial reports. 4th problem would have to be the fact that the chance of human error. Without a system in place to handle the orders, humans are going to be the ones who have to do the processing. So, human error are going to be an immanent problem. Human error itself can be considered an infinitely occurring liability for the company. The 5th and final problem would have to be the lack of progress reporting that the customers would need. For example, let" + [char]8217 + "s say that a customer buys a game on the 15th of December and we tell him/her that their game won" + [char]8217 + "t be until the 28th. As of this point, the customer is expected to receive their package on the 28th but suddenly something goes wrong. Let" + [char]8217 + "s say that the truck that the package is currently being transported in got into an accident which lead to the contents of the vehicle to be ruined. There" + [char]8217 + "s no way of us or the customer knowing of this until way later when the courier company itself calls us to let us know. Had we had a system that can accurately track the package, we" + [char]8217 + "d be able to know as the very minute the accident happens which would then give us enough time to let our customer know beforehand. No one gets left in the dark, everyone is informed, and everyone is happy. Next is the requirement phase. This part should be easy enough. First requirement would be the website. The website is needed by the customer to communicate with our system. It" + [char]8217 + "ll basically act as a window between the customer and the company."

$d.Content.Find.Execute(
    "which would " , $true, $false, $false, $false, $false, $true, 1, $false,
    ("which would " + $newBlock), 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Mark the five ordinal suffixes ("th") that follow digits as superscript,
#    matching Word's AutoFormat-as-you-type ordinal behaviour. Each is found
#    via a unique surrounding anchor, then the range is narrowed down to just
#    the "th" characters before the Superscript flag is applied.
# ---------------------------------------------------------------------------
$ordinals = @(
    @{ Anchor = "reports. 4th problem";    Offset = 10 },
    @{ Anchor = "The 5th and final";       Offset = 5 },
    @{ Anchor = "the 15th of December";    Offset = 6 },
    @{ Anchor = "until the 28th. As";      Offset = 12 },
    @{ Anchor = "the 28th but suddenly";   Offset = 6 }
)

foreach ($item in $ordinals) {
    $r = $d.Content
    $r.Find.Execute($item.Anchor) | Out-Null
    $base = $r.Start
    $r.Start = $base + $item.Offset
    $r.End = $base + $item.Offset + 2
    $r.Font.Superscript = $true
}

Write-Host "done"
